$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new Git branch/remote/conflict reference rows (49-68) ---
$ws.Range("A49").Value = 'git branch'
$ws.Range("B49").Value = 'shows all branches in a project'
$ws.Range("A50").Value = 'git checkout -b <branch name>'
$ws.Range("B50").Value = 'create a new branch'
$ws.Range("A51").Value = 'git diff main <branch name>'
$ws.Range("B51").Value = 'check difference between main and branch'
$ws.Range("A52").Value = 'git checkout <branch name>'
$ws.Range("B52").Value = 'switch to a branch'
$ws.Range("A53").Value = 'git merge <source> <destination>'
$ws.Range("B53").Value = 'merge branches e.g. git merge branch123 main'
$ws.Range("C49").Value = 'Branch'
$ws.Range("C50").Value = 'Branch'
$ws.Range("C51").Value = 'Branch'
$ws.Range("C52").Value = 'Branch'
$ws.Range("C53").Value = 'Branch'
$ws.Range("C54").Value = 'Conflicts'
$ws.Range("C55").Value = 'Conflicts'
$ws.Range("C56").Value = 'Conflicts'
$ws.Range("A54").Value = '<<<<<<< HEAD'
$ws.Range("B54").Value = 'Marks the start of the section with changes from your current branch (HEAD)'
$ws.Range("A55").Value = '''======='
$ws.Range("B55").Value = 'Divides the conflicting changes.'
$ws.Range("A56").Value = '>>>>>>> other-branch'
$ws.Range("B56").Value = 'Marks the end of the section with changes from the branch being merged (other-branch).'
$ws.Range("A57").Value = 'git init <folder name>'
$ws.Range("C57").Value = 'Repository'
$ws.Range("C58").Value = 'Repository'
$ws.Range("A58").Value = 'git init'
$ws.Range("B58").Value = 'Creates a Git repo for the current folder'
$ws.Range("B57").Value = 'Create a Git repo for the specified folder/project'
$ws.Range("A59").Value = 'git clone <repo path>'
$ws.Range("B59").Value = 'Clones the given repo'
$ws.Range("A60").Value = 'git clone <repo path> <repo name>'
$ws.Range("B60").Value = 'Clones the given repo with given name'
$ws.Range("C59").Value = 'Remote'
$ws.Range("C60").Value = 'Remote'
$ws.Range("A61").Value = 'git remote'
$ws.Range("B61").Value = 'lists name of remotes'
$ws.Range("A62").Value = 'git remote -v'
$ws.Range("B62").Value = 'returns remote url'
$ws.Range("B62").Interior.Color = 65535
$ws.Range("A63").Value = 'git remote add <name> <URL>'
$ws.Range("B63").Value = 'Renames ''main'' to the given name (Git automatically names remote as ''main'')'
$ws.Range("A64").Value = 'git fetch origin main'
$ws.Range("A65").Value = 'git fetch origin <repo name>'
$ws.Range("B64").Value = 'fetch from origin remote into local repo''s main branch'
$ws.Range("B65").Value = 'fetch from origin remote into local repo''s specific branch'
$ws.Range("A66").Value = 'git merge origin main'
$ws.Range("B66").Value = 'sync contents between remote and local main branch'
$ws.Range("B67").Value = 'pull from remote repo to main local branch'
$ws.Range("A68").Value = 'git  push <remote> <local branch>'
$ws.Range("A67").Value = 'git  pull <remote> <local branch>'
$ws.Range("B68").Value = 'push from local branch to remote repo'

# --- Update view: selection + scroll position ---
# Best-effort: scroll so row 54 is at the top of the view (matches the
# author's saved window state). Older/limited hosts may not persist this,
# so it's wrapped defensively and must not abort the rest of the script.
try { $excel.ActiveWindow.ScrollRow = 54 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("B69").Select()

# --- Page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1
